$wb = $excel.ActiveWorkbook

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1225
$ws.Range("I94").Value = 1225
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1225
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -774
$ws.Range("N94").ClearContents()

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1717.5209
$ws.Range("J138").Value = 2489.2856
$ws.Range("L138").Value = 7467.8568
$ws.Range("N138").Value = -17747.8568

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2287.6667
$ws.Range("I45").Value = 1244
$ws.Range("J45").Value = 4375
$ws.Range("K45").Value = 1244
$ws.Range("L45").Value = 4375
$ws.Range("M45").Value = -867
$ws.Range("N45").Value = -5129

# ARM row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 64988
$ws.Range("J128").Value = 64988
$ws.Range("L128").Value = 64988
$ws.Range("N128").Value = -74948

# ARM row 130
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 95335.664
$ws.Range("J130").Value = 95335.664
$ws.Range("L130").Value = 95335.664
$ws.Range("N130").Value = -105375.664

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2598.2354
$ws.Range("I132").Value = 2355.6
$ws.Range("J132").Value = 4418
$ws.Range("K132").Value = 7066.799999999999
$ws.Range("L132").Value = 13254
$ws.Range("M132").Value = -4536.799999999999
$ws.Range("N132").Value = -18314

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2391.3809
$ws.Range("I99").Value = 1593.8462
$ws.Range("K99").Value = 1593.8462
$ws.Range("M99").Value = -95.84619999999995

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2805.238
$ws.Range("I105").Value = 2393
$ws.Range("K105").Value = 2393
$ws.Range("M105").Value = -646

# BSM row 127
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H127").Value = 82781.5
$ws.Range("J127").Value = 82781.5
$ws.Range("L127").Value = 82781.5
$ws.Range("N127").Value = -92701.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6387.683
$ws.Range("I31").Value = 2556.077
$ws.Range("K31").Value = 2556.077
$ws.Range("M31").Value = -2261.077

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 6387.683
$ws.Range("I34").Value = 2556.077
$ws.Range("K34").Value = 2556.077
$ws.Range("M34").Value = -2354.077

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2401.8333
$ws.Range("I99").Value = 2082.2
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 2082.2
$ws.Range("L99").Value = 4000
$ws.Range("M99").Value = -584.1999999999998
$ws.Range("N99").Value = -6996

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2401.8333
$ws.Range("I126").Value = 2082.2
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 6246.599999999999
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -3776.599999999999
$ws.Range("N126").Value = -16940

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1652.25
$ws.Range("I134").Value = 1652.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4956.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2421.75
$ws.Range("N134").ClearContents()

# CUL row 7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 784.1429000000001
$ws.Range("I7").Value = 176.33333
$ws.Range("J7").Value = 949.9091
$ws.Range("K7").Value = 528.99999
$ws.Range("L7").Value = 2849.7273
$ws.Range("M7").Value = -416.99999
$ws.Range("N7").Value = -3073.7273

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2383.5557
$ws.Range("I33").Value = 2208.6
$ws.Range("J33").Value = 2602.25
$ws.Range("K33").Value = 13251.6
$ws.Range("L33").Value = 15613.5
$ws.Range("M33").Value = -12968.6
$ws.Range("N33").Value = -16179.5

# CUL row 47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1003.5
$ws.Range("I47").Value = 1003
$ws.Range("J47").Value = 1004
$ws.Range("K47").Value = 3009
$ws.Range("L47").Value = 3012
$ws.Range("M47").Value = -2578
$ws.Range("N47").Value = -3874

# CUL row 50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 523.2222
$ws.Range("I50").Value = 351.25
$ws.Range("J50").Value = 660.8
$ws.Range("K50").Value = 1053.75
$ws.Range("L50").Value = 1982.4
$ws.Range("M50").Value = -572.75
$ws.Range("N50").Value = -2944.4

# CUL row 53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 523.2222
$ws.Range("I53").Value = 351.25
$ws.Range("J53").Value = 660.8
$ws.Range("K53").Value = 1053.75
$ws.Range("L53").Value = 1982.4
$ws.Range("M53").Value = -572.75
$ws.Range("N53").Value = -2944.4

# CUL row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 639.5
$ws.Range("I116").Value = 394
$ws.Range("J116").Value = 885
$ws.Range("K116").Value = 1182
$ws.Range("L116").Value = 2655
$ws.Range("M116").Value = 2260
$ws.Range("N116").Value = -9539

# CUL row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 1579.6666
$ws.Range("I130").Value = 1579.6666
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 4738.9998
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 281.0002000000004
$ws.Range("N130").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1800
$ws.Range("J131").Value = 2500
$ws.Range("L131").Value = 7500
$ws.Range("N131").Value = -17580

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2509.5293
$ws.Range("I140").Value = 1151.091
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 3453.273
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 1726.727
$ws.Range("N140").Value = -25360

# GSM row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3193526
$ws.Range("J11").Value = 1090142
$ws.Range("L11").Value = 1090142
$ws.Range("N11").Value = -1090420

# GSM row 52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 47330.668
$ws.Range("J52").Value = 47246.25
$ws.Range("L52").Value = 47246.25
$ws.Range("N52").Value = -47764.25

# GSM row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# GSM row 131
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4721.091
$ws.Range("I132").Value = 4520.222
$ws.Range("J132").Value = 5625
$ws.Range("K132").Value = 13560.666
$ws.Range("L132").Value = 16875
$ws.Range("M132").Value = -11030.666
$ws.Range("N132").Value = -21935

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5280.3335
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# LTW row 11
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 12000
$ws.Range("J11").Value = 12000
$ws.Range("L11").Value = 12000
$ws.Range("N11").Value = -12280

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2701.7693
$ws.Range("J61").Value = 2578.125
$ws.Range("L61").Value = 2578.125
$ws.Range("N61").Value = -2982.125

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2701.7693
$ws.Range("J113").Value = 2578.125
$ws.Range("L113").Value = 2578.125
$ws.Range("N113").Value = -6918.125

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5280.3335
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3428.65
$ws.Range("I132").Value = 3197.8462
$ws.Range("J132").Value = 3857.2856
$ws.Range("K132").Value = 9593.5386
$ws.Range("L132").Value = 11571.8568
$ws.Range("M132").Value = -7063.5386
$ws.Range("N132").Value = -16631.8568

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2628
$ws.Range("I136").Value = 2375.0908
$ws.Range("J136").Value = 3184.4
$ws.Range("K136").Value = 7125.2724
$ws.Range("L136").Value = 9553.200000000001
$ws.Range("M136").Value = -4575.2724
$ws.Range("N136").Value = -14653.2

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 295997.8
$ws.Range("J140").Value = 295997.8
$ws.Range("L140").Value = 295997.8
$ws.Range("N140").Value = -306357.8

# WVR row 38
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()

# WVR row 49
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 30053
$ws.Range("J49").Value = 30050
$ws.Range("L49").Value = 30050
$ws.Range("N49").Value = -30510

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1374.75
$ws.Range("I100").Value = 876
$ws.Range("J100").Value = 1873.5
$ws.Range("K100").Value = 1752
$ws.Range("L100").Value = 3747
$ws.Range("M100").Value = -1211
$ws.Range("N100").Value = -4829

# WVR row 124
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 99999
$ws.Range("J124").Value = 99999
$ws.Range("L124").Value = 99999
$ws.Range("N124").Value = -109819
